$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Preserve text storage (these columns hold text-formatted numbers/percentages)
# by forcing Text number format before assigning numeric-looking strings,
# matching the workbook pre-existing t="inlineStr" (text) cell type.
$ws.Range("D2:D51").NumberFormat = "@"
$ws.Range("E2:E51").NumberFormat = "@"
$ws.Range("G2:G51").NumberFormat = "@"

$ws.Cells.Item(2, 4).Value = "261.33"
$ws.Cells.Item(2, 5).Value = "0.75%"
$ws.Cells.Item(2, 7).Value = "23"
$ws.Cells.Item(3, 4).Value = "27.14"
$ws.Cells.Item(3, 5).Value = "1.06%"
$ws.Cells.Item(3, 7).Value = "23"
$ws.Cells.Item(4, 4).Value = "4.708"
$ws.Cells.Item(4, 5).Value = "0.77%"
$ws.Cells.Item(4, 7).Value = "23"
$ws.Cells.Item(5, 4).Value = "0.06207"
$ws.Cells.Item(5, 5).Value = "2.46%"
$ws.Cells.Item(5, 7).Value = "23"
$ws.Cells.Item(6, 4).Value = "6.732"
$ws.Cells.Item(6, 5).Value = "0.65%"
$ws.Cells.Item(6, 7).Value = "23"
$ws.Cells.Item(7, 4).Value = "0.8509"
$ws.Cells.Item(7, 5).Value = "-1.25%"
$ws.Cells.Item(7, 7).Value = "23"
$ws.Cells.Item(8, 4).Value = "0.9071"
$ws.Cells.Item(8, 5).Value = "-1.62%"
$ws.Cells.Item(8, 7).Value = "23"
$ws.Cells.Item(9, 4).Value = "0.1403"
$ws.Cells.Item(9, 5).Value = "0.38%"
$ws.Cells.Item(9, 7).Value = "23"
$ws.Cells.Item(10, 4).Value = "0.04776"
$ws.Cells.Item(10, 5).Value = "-7.43%"
$ws.Cells.Item(10, 7).Value = "23"
$ws.Cells.Item(11, 4).Value = "0.07099"
$ws.Cells.Item(11, 5).Value = "0.07%"
$ws.Cells.Item(11, 7).Value = "23"
$ws.Cells.Item(12, 4).Value = "0.03165"
$ws.Cells.Item(12, 5).Value = "1.34%"
$ws.Cells.Item(12, 7).Value = "23"
$ws.Cells.Item(13, 4).Value = "0.09056"
$ws.Cells.Item(13, 5).Value = "-0.81%"
$ws.Cells.Item(13, 7).Value = "23"
$ws.Cells.Item(14, 4).Value = "0.001539"
$ws.Cells.Item(14, 5).Value = "-0.40%"
$ws.Cells.Item(14, 7).Value = "23"
$ws.Cells.Item(15, 4).Value = "0.0006152"
$ws.Cells.Item(15, 5).Value = "1.45%"
$ws.Cells.Item(15, 7).Value = "23"
$ws.Cells.Item(16, 4).Value = "0.006134"
$ws.Cells.Item(16, 5).Value = "2.14%"
$ws.Cells.Item(16, 7).Value = "23"
$ws.Cells.Item(17, 4).Value = "3.468"
$ws.Cells.Item(17, 5).Value = "-0.45%"
$ws.Cells.Item(17, 7).Value = "23"
$ws.Cells.Item(18, 4).Value = "3.170"
$ws.Cells.Item(18, 5).Value = "0.00%"
$ws.Cells.Item(18, 7).Value = "23"
$ws.Cells.Item(19, 4).Value = "2.177"
$ws.Cells.Item(19, 5).Value = "-0.35%"
$ws.Cells.Item(19, 7).Value = "23"
$ws.Cells.Item(20, 7).Value = "23"
$ws.Cells.Item(21, 5).Value = "-1.26%"
$ws.Cells.Item(21, 7).Value = "23"
$ws.Cells.Item(22, 4).Value = "4.118"
$ws.Cells.Item(22, 5).Value = "0.78%"
$ws.Cells.Item(22, 7).Value = "23"
$ws.Cells.Item(23, 4).Value = "0.04245"
$ws.Cells.Item(23, 5).Value = "0.18%"
$ws.Cells.Item(23, 7).Value = "23"
$ws.Cells.Item(24, 4).Value = "0.001217"
$ws.Cells.Item(24, 5).Value = "0.02%"
$ws.Cells.Item(24, 7).Value = "23"
$ws.Cells.Item(25, 4).Value = "0.004117"
$ws.Cells.Item(25, 5).Value = "2.41%"
$ws.Cells.Item(25, 7).Value = "23"
$ws.Cells.Item(26, 5).Value = "0.15%"
$ws.Cells.Item(26, 7).Value = "23"
$ws.Cells.Item(27, 7).Value = "23"
$ws.Cells.Item(28, 7).Value = "23"
$ws.Cells.Item(29, 7).Value = "23"
$ws.Cells.Item(30, 7).Value = "23"
$ws.Cells.Item(31, 7).Value = "23"
$ws.Cells.Item(32, 7).Value = "23"
$ws.Cells.Item(33, 7).Value = "23"
$ws.Cells.Item(34, 7).Value = "23"
$ws.Cells.Item(35, 7).Value = "23"
$ws.Cells.Item(36, 7).Value = "23"
$ws.Cells.Item(37, 7).Value = "23"
$ws.Cells.Item(38, 7).Value = "23"
$ws.Cells.Item(39, 7).Value = "23"
$ws.Cells.Item(40, 4).Value = "0.03903"
$ws.Cells.Item(40, 5).Value = "0.87%"
$ws.Cells.Item(40, 7).Value = "23"
$ws.Cells.Item(41, 4).Value = "0.1112"
$ws.Cells.Item(41, 5).Value = "-0.41%"
$ws.Cells.Item(41, 7).Value = "23"
$ws.Cells.Item(42, 4).Value = "0.004136"
$ws.Cells.Item(42, 5).Value = "1.70%"
$ws.Cells.Item(42, 7).Value = "23"
$ws.Cells.Item(43, 5).Value = "-0.69%"
$ws.Cells.Item(43, 7).Value = "23"
$ws.Cells.Item(44, 4).Value = "0.01341"
$ws.Cells.Item(44, 5).Value = "-10.12%"
$ws.Cells.Item(44, 7).Value = "23"
$ws.Cells.Item(45, 5).Value = "-0.24%"
$ws.Cells.Item(45, 7).Value = "23"
$ws.Cells.Item(46, 5).Value = "0.14%"
$ws.Cells.Item(46, 7).Value = "23"
$ws.Cells.Item(47, 5).Value = "-34.14%"
$ws.Cells.Item(47, 7).Value = "23"
$ws.Cells.Item(48, 4).Value = "0.05782"
$ws.Cells.Item(48, 5).Value = "-57.26%"
$ws.Cells.Item(48, 7).Value = "23"
$ws.Cells.Item(49, 5).Value = "0.14%"
$ws.Cells.Item(49, 7).Value = "23"
$ws.Cells.Item(50, 5).Value = "0.14%"
$ws.Cells.Item(50, 7).Value = "23"
$ws.Cells.Item(51, 7).Value = "23"
